$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column BB is added after BA (column 53 -> 54), replicating column BA's
# date-header formatting for row 1 and carrying forward the YoY forecast
# series, with rows 19-21 (the newest forecast periods) updated to reflect
# the re-run ("Included EQUIPMENT eval ... allowed for multiple archive excels").

# Row 1: new quarter header date, formatted like the rest of the header row.
$ws.Range("BA1").Copy() | Out-Null
$ws.Range("BB1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("BB1").Value = 45986

# Rows 3-18: same value as column BA (forecast held flat / unchanged).
$sameValueRows = 3..18
foreach ($r in $sameValueRows) {
    $baCell = $ws.Cells.Item($r, 53)
    $bbCell = $ws.Cells.Item($r, 54)
    $bbCell.Value = $baCell.Value2
}

# Rows 19-21: updated forecast values from the new run.
$ws.Cells.Item(19, 54).Value = 2.560577522109297
$ws.Cells.Item(20, 54).Value = 1.325305149734723
$ws.Cells.Item(21, 54).Value = 0.1333065884001616

Write-Host "Applied BB column update"
